$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the viewport one column right so column V becomes the leftmost
# visible column (topLeftCell "U1" -> "V1").
$excel.ActiveWindow.ScrollColumn = 22
$excel.ActiveWindow.ScrollRow = 1

# AG8 gets its own standalone formula.
$ws.Range("AG8").Formula = '=TEXT(TODAY(), "dd-mmm-yyyy")'

# AG9:AG12 are filled with the same formula as one shared-formula group.
$ws.Range("AG9:AG12").Formula = '=TEXT(TODAY(), "dd-mmm-yyyy")'

# Move the selection to AF16.
$ws.Range("AF16").Select() | Out-Null
